# Applies "Atualização arquivo de falas":
#  - Slide 2 paragraph ("O nosso negocio..."): the sentence had been typed
#    across six separate runs; collapse them into one run with the full
#    sentence as its text (formatting sz=24/szCs=24 unchanged).
#  - Slide 2 paragraph ("Atualmente uma das grandes dificuldades..."):
#    the description of what matters in the silos is reworded from
#    "temperatura, umidade e luminosidade" to "temperatura e umidade",
#    with the text re-split into runs around the existing _GoBack
#    bookmark, and the closing curly quote folded into the final run.

$d = $word.ActiveDocument

$p5 = $d.Paragraphs.Item(5)
$xml5 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="009F7BF6" w:rsidRDefault="004E7A7F" w:rsidP="000007DA"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">“O nosso negócio, a Coffee Tech - Auditoria de armazenagem computadorizada, é uma empresa que surgiu para atender uma dor do mercado no segmento de produção de café, mais especificadamente no setor de armazenagem. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $p5.Range.InsertXML($xml5)

$p6 = $d.Paragraphs.Item(6)
$xml6 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="000007DA" w:rsidRPr="004E7A7F" w:rsidRDefault="004E7A7F" w:rsidP="000007DA"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="004E7A7F"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Atualmente uma das grandes dificuldades das empresas </w:t></w:r><w:r w:rsidR="007F4356"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>desse</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> setor, é se certificarem que de fato os níveis de temperatura</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> e</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> umidade</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>nos silos de armazenamento de grãos estão corretos. E quem vai falar um pouco mais sobre esse problema, é o meu sócio Diego.”</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $p6.Range.InsertXML($xml6)
